# ============================================================================
# 688107-安路科技.xlsx -- add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" sheet right after "总计", before "2022-Q3" (so it
#    becomes the 2nd tab; the existing quarter sheets just shift right).
# 2) Prepend a 2022-Q4 summary row into "总计", pushing the other rows down.
# 3) Populate the new sheet with the 2022-Q4 fund-holdings table.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet -- shift rows 2-5 down to 3-6,
#         keep the sequential index column (A) as-is, then write the new
#         2022-Q4 row into row 2.
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B2:D5").Copy($ws1.Range("B3:D6"))

# A6 is a brand-new row; clone A5's formatting (bold / centred / bordered)
# onto it before writing its index value.
$ws1.Range("A5").Copy($ws1.Range("A6"))
$ws1.Cells.Item(6,1).Value = 4

$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,3).Value = 30
$ws1.Cells.Item(2,4).Value = 10.07

# ----------------------------------------------------------------------------
# Step 2: insert the new "2022-Q4" worksheet as the 2nd tab.
# ----------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q4"

# Borrow header-row (B1:H1) + index-column (A2:A31) formatting from the
# neighbouring quarter sheet (now shifted to tab #3) so the new sheet matches
# the established look of the other quarter tabs.
$sheetQ3 = $wb.Worksheets.Item(3)
$sheetQ3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$sheetQ3.Range("A2").Copy($newSheet.Range("A2:A31"))

# ----------------------------------------------------------------------------
# Step 3: fill the 30 fund-holding rows (row 2 .. row 31).
# Columns B,D,E,F,G hold numeric-looking text in the source data; force them
# to text storage (t="inlineStr"), same as the other quarter sheets, then drop
# back to the default style. Row 31 col G (011702 / 广发睿享稳健增利混合C) is the
# sole exception: its market-value cell is a genuine number 0 in the source.
# ----------------------------------------------------------------------------
$newSheet.Range("B2:B31").NumberFormat = "@"
$newSheet.Range("D2:F31").NumberFormat = "@"
$newSheet.Range("G2:G30").NumberFormat = "@"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "519674"
$newSheet.Cells.Item(2,3).Value = "银河创新成长混合A"
$newSheet.Cells.Item(2,4).Value = "145.89"
$newSheet.Cells.Item(2,5).Value = "92.48"
$newSheet.Cells.Item(2,6).Value = "3.41"
$newSheet.Cells.Item(2,7).Value = "4.9748"
$newSheet.Cells.Item(2,8).Value = 10
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "007872"
$newSheet.Cells.Item(3,3).Value = "金信稳健策略灵活配置混合"
$newSheet.Cells.Item(3,4).Value = "26.03"
$newSheet.Cells.Item(3,5).Value = "93.83"
$newSheet.Cells.Item(3,6).Value = "8.41"
$newSheet.Cells.Item(3,7).Value = "2.1891"
$newSheet.Cells.Item(3,8).Value = 3
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "014143"
$newSheet.Cells.Item(4,3).Value = "银河创新成长混合C"
$newSheet.Cells.Item(4,4).Value = "22.12"
$newSheet.Cells.Item(4,5).Value = "92.48"
$newSheet.Cells.Item(4,6).Value = "3.41"
$newSheet.Cells.Item(4,7).Value = "0.7543"
$newSheet.Cells.Item(4,8).Value = 10
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "007490"
$newSheet.Cells.Item(5,3).Value = "南方信息创新混合A"
$newSheet.Cells.Item(5,4).Value = "15.05"
$newSheet.Cells.Item(5,5).Value = "91.51"
$newSheet.Cells.Item(5,6).Value = "3.85"
$newSheet.Cells.Item(5,7).Value = "0.5794"
$newSheet.Cells.Item(5,8).Value = 10
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "002810"
$newSheet.Cells.Item(6,3).Value = "金信转型创新成长灵活配置混合"
$newSheet.Cells.Item(6,4).Value = "3.84"
$newSheet.Cells.Item(6,5).Value = "89.18"
$newSheet.Cells.Item(6,6).Value = "8.45"
$newSheet.Cells.Item(6,7).Value = "0.3245"
$newSheet.Cells.Item(6,8).Value = 1
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "002256"
$newSheet.Cells.Item(7,3).Value = "金信行业优选灵活配置混合"
$newSheet.Cells.Item(7,4).Value = "3.56"
$newSheet.Cells.Item(7,5).Value = "93.94"
$newSheet.Cells.Item(7,6).Value = "8.45"
$newSheet.Cells.Item(7,7).Value = "0.3008"
$newSheet.Cells.Item(7,8).Value = 3
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "660015"
$newSheet.Cells.Item(8,3).Value = "农银行业轮动混合A"
$newSheet.Cells.Item(8,4).Value = "6.67"
$newSheet.Cells.Item(8,5).Value = "86.12"
$newSheet.Cells.Item(8,6).Value = "2.68"
$newSheet.Cells.Item(8,7).Value = "0.1788"
$newSheet.Cells.Item(8,8).Value = 6
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "014185"
$newSheet.Cells.Item(9,3).Value = "招商专精特新股票A"
$newSheet.Cells.Item(9,4).Value = "3.30"
$newSheet.Cells.Item(9,5).Value = "87.72"
$newSheet.Cells.Item(9,6).Value = "4.69"
$newSheet.Cells.Item(9,7).Value = "0.1548"
$newSheet.Cells.Item(9,8).Value = 8
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "630010"
$newSheet.Cells.Item(10,3).Value = "华商价值精选混合"
$newSheet.Cells.Item(10,4).Value = "4.37"
$newSheet.Cells.Item(10,5).Value = "87.25"
$newSheet.Cells.Item(10,6).Value = "3.10"
$newSheet.Cells.Item(10,7).Value = "0.1355"
$newSheet.Cells.Item(10,8).Value = 8
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = "014186"
$newSheet.Cells.Item(11,3).Value = "招商专精特新股票C"
$newSheet.Cells.Item(11,4).Value = "2.50"
$newSheet.Cells.Item(11,5).Value = "87.72"
$newSheet.Cells.Item(11,6).Value = "4.69"
$newSheet.Cells.Item(11,7).Value = "0.1172"
$newSheet.Cells.Item(11,8).Value = 8
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = "007491"
$newSheet.Cells.Item(12,3).Value = "南方信息创新混合C"
$newSheet.Cells.Item(12,4).Value = "2.40"
$newSheet.Cells.Item(12,5).Value = "91.51"
$newSheet.Cells.Item(12,6).Value = "3.85"
$newSheet.Cells.Item(12,7).Value = "0.0924"
$newSheet.Cells.Item(12,8).Value = 10
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = "501076"
$newSheet.Cells.Item(13,3).Value = "鹏华创新动力混合（LOF）"
$newSheet.Cells.Item(13,4).Value = "5.11"
$newSheet.Cells.Item(13,5).Value = "51.28"
$newSheet.Cells.Item(13,6).Value = "1.27"
$newSheet.Cells.Item(13,7).Value = "0.0649"
$newSheet.Cells.Item(13,8).Value = 7
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = "005117"
$newSheet.Cells.Item(14,3).Value = "金信价值精选灵活配置混合A"
$newSheet.Cells.Item(14,4).Value = "0.76"
$newSheet.Cells.Item(14,5).Value = "92.96"
$newSheet.Cells.Item(14,6).Value = "3.68"
$newSheet.Cells.Item(14,7).Value = "0.0280"
$newSheet.Cells.Item(14,8).Value = 4
$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = "630006"
$newSheet.Cells.Item(15,3).Value = "华商产业升级混合"
$newSheet.Cells.Item(15,4).Value = "0.86"
$newSheet.Cells.Item(15,5).Value = "88.65"
$newSheet.Cells.Item(15,6).Value = "3.17"
$newSheet.Cells.Item(15,7).Value = "0.0273"
$newSheet.Cells.Item(15,8).Value = 8
$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = "012696"
$newSheet.Cells.Item(16,3).Value = "同泰数字经济主题股票A"
$newSheet.Cells.Item(16,4).Value = "0.87"
$newSheet.Cells.Item(16,5).Value = "93.38"
$newSheet.Cells.Item(16,6).Value = "3.01"
$newSheet.Cells.Item(16,7).Value = "0.0262"
$newSheet.Cells.Item(16,8).Value = 9
$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).Value = "007251"
$newSheet.Cells.Item(17,3).Value = "广发睿享稳健增利混合A"
$newSheet.Cells.Item(17,4).Value = "1.26"
$newSheet.Cells.Item(17,5).Value = "39.64"
$newSheet.Cells.Item(17,6).Value = "1.72"
$newSheet.Cells.Item(17,7).Value = "0.0217"
$newSheet.Cells.Item(17,8).Value = 10
$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).Value = "002292"
$newSheet.Cells.Item(18,3).Value = "诺安益鑫灵活配置混合A"
$newSheet.Cells.Item(18,4).Value = "0.37"
$newSheet.Cells.Item(18,5).Value = "69.58"
$newSheet.Cells.Item(18,6).Value = "5.38"
$newSheet.Cells.Item(18,7).Value = "0.0199"
$newSheet.Cells.Item(18,8).Value = 2
$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).Value = "015919"
$newSheet.Cells.Item(19,3).Value = "申万菱信专精特新主题混合A"
$newSheet.Cells.Item(19,4).Value = "0.39"
$newSheet.Cells.Item(19,5).Value = "48.55"
$newSheet.Cells.Item(19,6).Value = "4.10"
$newSheet.Cells.Item(19,7).Value = "0.0160"
$newSheet.Cells.Item(19,8).Value = 4
$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).Value = "015850"
$newSheet.Cells.Item(20,3).Value = "农银行业轮动混合C"
$newSheet.Cells.Item(20,4).Value = "0.52"
$newSheet.Cells.Item(20,5).Value = "86.12"
$newSheet.Cells.Item(20,6).Value = "2.68"
$newSheet.Cells.Item(20,7).Value = "0.0139"
$newSheet.Cells.Item(20,8).Value = 6
$newSheet.Cells.Item(21,1).Value = 19
$newSheet.Cells.Item(21,2).Value = "001226"
$newSheet.Cells.Item(21,3).Value = "中邮稳健添利灵活配置混合"
$newSheet.Cells.Item(21,4).Value = "0.40"
$newSheet.Cells.Item(21,5).Value = "93.24"
$newSheet.Cells.Item(21,6).Value = "2.91"
$newSheet.Cells.Item(21,7).Value = "0.0116"
$newSheet.Cells.Item(21,8).Value = 5
$newSheet.Cells.Item(22,1).Value = 20
$newSheet.Cells.Item(22,2).Value = "003238"
$newSheet.Cells.Item(22,3).Value = "新华外延增长主题灵活配置混合"
$newSheet.Cells.Item(22,4).Value = "0.51"
$newSheet.Cells.Item(22,5).Value = "85.17"
$newSheet.Cells.Item(22,6).Value = "2.23"
$newSheet.Cells.Item(22,7).Value = "0.0114"
$newSheet.Cells.Item(22,8).Value = 7
$newSheet.Cells.Item(23,1).Value = 21
$newSheet.Cells.Item(23,2).Value = "012697"
$newSheet.Cells.Item(23,3).Value = "同泰数字经济主题股票C"
$newSheet.Cells.Item(23,4).Value = "0.34"
$newSheet.Cells.Item(23,5).Value = "93.38"
$newSheet.Cells.Item(23,6).Value = "3.01"
$newSheet.Cells.Item(23,7).Value = "0.0102"
$newSheet.Cells.Item(23,8).Value = 9
$newSheet.Cells.Item(24,1).Value = 22
$newSheet.Cells.Item(24,2).Value = "519172"
$newSheet.Cells.Item(24,3).Value = "浦银安盛睿智精选灵活配置混合A"
$newSheet.Cells.Item(24,4).Value = "0.21"
$newSheet.Cells.Item(24,5).Value = "89.04"
$newSheet.Cells.Item(24,6).Value = "3.16"
$newSheet.Cells.Item(24,7).Value = "0.0066"
$newSheet.Cells.Item(24,8).Value = 4
$newSheet.Cells.Item(25,1).Value = 23
$newSheet.Cells.Item(25,2).Value = "519173"
$newSheet.Cells.Item(25,3).Value = "浦银安盛睿智精选灵活配置混合C"
$newSheet.Cells.Item(25,4).Value = "0.13"
$newSheet.Cells.Item(25,5).Value = "89.04"
$newSheet.Cells.Item(25,6).Value = "3.16"
$newSheet.Cells.Item(25,7).Value = "0.0041"
$newSheet.Cells.Item(25,8).Value = 4
$newSheet.Cells.Item(26,1).Value = 24
$newSheet.Cells.Item(26,2).Value = "005118"
$newSheet.Cells.Item(26,3).Value = "金信价值精选灵活配置混合C"
$newSheet.Cells.Item(26,4).Value = "0.06"
$newSheet.Cells.Item(26,5).Value = "92.96"
$newSheet.Cells.Item(26,6).Value = "3.68"
$newSheet.Cells.Item(26,7).Value = "0.0022"
$newSheet.Cells.Item(26,8).Value = 4
$newSheet.Cells.Item(27,1).Value = 25
$newSheet.Cells.Item(27,2).Value = "014550"
$newSheet.Cells.Item(27,3).Value = "诺安益鑫灵活配置混合C"
$newSheet.Cells.Item(27,4).Value = "0.02"
$newSheet.Cells.Item(27,5).Value = "69.58"
$newSheet.Cells.Item(27,6).Value = "5.38"
$newSheet.Cells.Item(27,7).Value = "0.0011"
$newSheet.Cells.Item(27,8).Value = 2
$newSheet.Cells.Item(28,1).Value = 26
$newSheet.Cells.Item(28,2).Value = "015920"
$newSheet.Cells.Item(28,3).Value = "申万菱信专精特新主题混合C"
$newSheet.Cells.Item(28,4).Value = "0.02"
$newSheet.Cells.Item(28,5).Value = "48.55"
$newSheet.Cells.Item(28,6).Value = "4.10"
$newSheet.Cells.Item(28,7).Value = "0.0008"
$newSheet.Cells.Item(28,8).Value = 4
$newSheet.Cells.Item(29,1).Value = 27
$newSheet.Cells.Item(29,2).Value = "166108"
$newSheet.Cells.Item(29,3).Value = "信澳量化多因子混合（LOF）C"
$newSheet.Cells.Item(29,4).Value = "0.07"
$newSheet.Cells.Item(29,5).Value = "34.40"
$newSheet.Cells.Item(29,6).Value = "1.06"
$newSheet.Cells.Item(29,7).Value = "0.0007"
$newSheet.Cells.Item(29,8).Value = 1
$newSheet.Cells.Item(30,1).Value = 28
$newSheet.Cells.Item(30,2).Value = "166107"
$newSheet.Cells.Item(30,3).Value = "信澳量化多因子混合（LOF）A"
$newSheet.Cells.Item(30,4).Value = "0.05"
$newSheet.Cells.Item(30,5).Value = "34.40"
$newSheet.Cells.Item(30,6).Value = "1.06"
$newSheet.Cells.Item(30,7).Value = "0.0005"
$newSheet.Cells.Item(30,8).Value = 1
$newSheet.Cells.Item(31,1).Value = 29
$newSheet.Cells.Item(31,2).Value = "011702"
$newSheet.Cells.Item(31,3).Value = "广发睿享稳健增利混合C"
$newSheet.Cells.Item(31,4).Value = "0.00"
$newSheet.Cells.Item(31,5).Value = "39.64"
$newSheet.Cells.Item(31,6).Value = "1.72"
$newSheet.Cells.Item(31,7).Value = 0
$newSheet.Cells.Item(31,8).Value = 10

# Reset the forced text-format cells back to the default (General) style now
# that the values are committed as text -- matches the target, which has no
# explicit "s" attribute on these data cells.
$newSheet.Range("B2:B31").Style = "Normal"
$newSheet.Range("D2:F31").Style = "Normal"
$newSheet.Range("G2:G30").Style = "Normal"

# ----------------------------------------------------------------------------
# Keep "总计" as the active/selected tab (matches the original workbook).
# ----------------------------------------------------------------------------
$ws1.Activate()
